$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final sorted dataset (ascending by column A), now spanning rows 2-10
$data = @(
    @(-0.02125, 0.30512),
    @(-0.01813, 0.39656),
    @(-0.011,   0.59931),
    @(-0.00875, 0.67123),
    @(-0.00625, 0.73175),
    @(-0.00562, 0.7526),
    @(-0.005,   0.76936),
    @(-0.00375, 0.80162),
    @(0,        0.90934)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Update the selected cell/range to match the saved view state
$ws.Range("B15").Select()
